$d = $word.ActiveDocument

# Anchor on the paragraph right before the block to remove:
# "...LOQ4044: Introducao a Engenharia da Qualidade (Requisito fraco)"
$anchorRange = $d.Content.Duplicate
$anchorRange.Find.Execute("LOQ4044", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0) | Out-Null
$anchorRange.Expand(4) | Out-Null  # wdParagraph
$anchorEnd = $anchorRange.End      # end of that paragraph, incl. its paragraph mark

# Anchor on the last paragraph of the block to remove: the copyright notice paragraph
# "(c) 2020 . Contact: ... Powered by Jekyll and Github pages. ..."
$endRange = $d.Content.Duplicate
$endRange.Find.Execute("Powered by Jekyll", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$endRange.Expand(4) | Out-Null     # wdParagraph
$blockEnd = $endRange.End          # end of that paragraph, incl. its paragraph mark

# Remove the blank paragraph, the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph,
# and the copyright paragraph -- i.e. everything between the two anchors above. The blank
# paragraph that precedes the trailing page-break paragraph is left untouched.
$rng = $d.Range($anchorEnd, $blockEnd)
$rng.Delete()
